# Updated cryptos list on Thu May 30 17:43:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "599.00", "171.96") are not coerced into floating point
    # numbers, then restore the default "Normal" style so no stray
    # number-format styling is left behind on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.295.47"
$ws.Range("E2").Value = "  +2.52%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.818.60"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
Set-TextValue "D5" "599.00"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6 - Solana
Set-TextValue "D6" "171.96"
$ws.Range("E6").Value = "  +1.00%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.817.68"
$ws.Range("E7").Value = "  +1.33%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.08%  "

# Row 9 - XRP
Set-TextValue "D9" "0.526"
$ws.Range("E9").Value = "  -0.49%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.82%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.55"
$ws.Range("E11").Value = "  +0.96%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.454"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000268"
$ws.Range("E13").Value = "  -3.63%  "

# Row 14 - Avalanche
Set-TextValue "D14" "37.05"
$ws.Range("E14").Value = "  +0.76%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.458.53"
$ws.Range("E15").Value = "  +1.24%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.815.96"
$ws.Range("E16").Value = "  +1.38%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "69.186.22"
$ws.Range("E17").Value = "  +2.30%  "

# Row 18 - Chainlink
Set-TextValue "D18" "18.37"
$ws.Range("E18").Value = "  -2.75%  "

# Row 19 - Polkadot
Set-TextValue "D19" "7.13"
$ws.Range("E19").Value = "  -1.70%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.09%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.17"
$ws.Range("E21").Value = "  +5.12%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "473.95"
$ws.Range("E22").Value = "  +0.81%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.712"
$ws.Range("E23").Value = "  -1.51%  "

# Row 24 - now PEPE (was Litecoin)
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D24" "0.0000149"
$ws.Range("E24").Value = "  +0.66%  "

# Row 25 - now Litecoin (was PEPE)
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "85.10"
$ws.Range("E25").Value = "  +1.33%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +0.72%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.28"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -1.17%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.08%  "

# Row 30 - WrappedeETH
Set-TextValue "D30" "3.967.22"
$ws.Range("E30").Value = "  +1.37%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.76%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "7.52"
$ws.Range("E32").Value = "  -2.74%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  +0.41%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "30.47"
$ws.Range("E34").Value = "  -0.16%  "

# Row 35 - Aptos
Set-TextValue "D35" "9.44"
$ws.Range("E35").Value = "  +2.83%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.772.43"
$ws.Range("E37").Value = "  +0.98%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -2.71%  "

# Row 39 - dogwifhat
Set-TextValue "D39" "3.61"
$ws.Range("E39").Value = "  -6.70%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +1.85%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  +0.84%  "

# Row 42 - Filecoin
Set-TextValue "D42" "5.92"
$ws.Range("E42").Value = "  +0.32%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  -0.08%  "

# Row 44 - TheGraph
$ws.Range("E44").Value = "  -1.16%  "

# Row 46 - now Stacks (was Arweave)
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "1.99"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47 - now Arweave (was Stacks)
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D47" "44.21"
$ws.Range("E47").Value = "  +13.08%  "

# Row 48 - Cosmos
Set-TextValue "D48" "8.67"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49 - now Bittensor (was OKB)
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D49" "407.19"
$ws.Range("E49").Value = "  +2.01%  "

# Row 50 - now OKB (was Bittensor)
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D50" "46.42"
$ws.Range("E50").Value = "  +1.11%  "

# Row 51 - Monero
Set-TextValue "D51" "146.10"
